$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Update the existing Mon..Fri time values first (columns still B..F at this
# point) before the new Sat/Sun columns are inserted.
# ---------------------------------------------------------------------------

# Row 2 "Time In": Tue/Wed/Thu share "10:00 AM", Fri becomes "9:00 AM", and
# Mon gets the no-space "10:00AM" variant.
$ws.Range("C2").Value = "10:00 AM"
$ws.Range("F2").Value = "9:00 AM"

# Row 5 "Time Out": Tue/Wed/Thu share "6:00 PM", Fri becomes "5:00 PM", and
# Mon gets the no-space "6:00PM" variant.
$ws.Range("C5").Value = "6:00 PM"
$ws.Range("F5").Value = "5:00 PM"

$ws.Range("B2").Value = "10:00AM"
$ws.Range("B5").Value = "6:00PM"

$ws.Range("D2").Value = "10:00 AM"
$ws.Range("E2").Value = "10:00 AM"
$ws.Range("D5").Value = "6:00 PM"
$ws.Range("E5").Value = "6:00 PM"

# Rows 3 & 4 "Meal Break Out" / "Meal Break In": blank these out.
$ws.Range("B3:F3").ClearContents()
$ws.Range("B4:F4").ClearContents()

# ---------------------------------------------------------------------------
# Insert two new columns before column B. This shifts the existing Mon..Fri
# columns (B..F) over to D..H, making room for new Sat/Sun columns.
# ---------------------------------------------------------------------------
$ws.Range("B1:C1").EntireColumn.Insert()

# The column insert blindly copies the formatting of the (new, now-blank)
# B/C columns down into every row, including row 6 and the now-empty B/C
# cells in the data rows. Clean that collateral damage up so only the cells
# that should exist in the final sheet remain.
$ws.Range("B2:C5").Style = "Normal"
$ws.Range("B2:C5").ClearContents()

$ws.Range("C6:D6").Style = "Normal"
$ws.Range("C6:D6").ClearContents()
$ws.Range("B6").Style = "Normal"
$ws.Range("B6").NumberFormat = "h:mm"

# ---------------------------------------------------------------------------
# New header cells: Sat / Sun, styled like the other day headers (copy the
# plain, non-wrapped header look from D1 "Mon") but with a red fill (re-using
# the existing white header font already in the file).
# ---------------------------------------------------------------------------
$ws.Range("D1").Copy()
$ws.Range("B1").PasteSpecial(-4122)
$ws.Range("B1").Value = "Sat"
$ws.Range("B1").Interior.Color = 255

$ws.Range("D1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$ws.Range("C1").Value = "Sun"
$ws.Range("C1").Interior.Color = 255

# ---------------------------------------------------------------------------
# Rows 1-5 got a bit taller in the revised layout.
# ---------------------------------------------------------------------------
$ws.Range("A1:H5").RowHeight = 16

# ---------------------------------------------------------------------------
# Selection cursor, mirroring the author's final cursor position.
# ---------------------------------------------------------------------------
$ws.Range("P7").Select()
